$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17 (existing rows 17-54 shift down to 18-55).
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with a new weekly price observation.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44645
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100108
$ws.Cells.Item(17, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(17, 9).Value = 100108007
$ws.Cells.Item(17, 10).Value = "Coco"
$ws.Cells.Item(17, 11).Value = "Sin especificar"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 10
$ws.Cells.Item(17, 14).Value = 28000
$ws.Cells.Item(17, 15).Value = 28000
$ws.Cells.Item(17, 16).Value = 28000
$ws.Cells.Item(17, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(17, 18).Value = "Perú"
$ws.Cells.Item(17, 19).Value = 1400
$ws.Cells.Item(17, 20).Value = 20
